# refactor: email field's biz in user module
# - make email optional (was required)
# - add banEmailUpdate / banThirdParty columns
# - shuffle header/description columns accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: field names (header row) -----------------------------------
$ws.Range("E1").Value = "phone"
$ws.Range("F1").Value = "gender"
$ws.Range("G1").Value = "studentId"
$ws.Range("H1").Value = "banThirdParty"
$ws.Range("I1").Value = "banEmailUpdate"

# --- Row 2: field descriptions (sample / hint row) ----------------------
$ws.Range("C2").Value = "选填，邮箱"
$ws.Range("E2").Value = "选填，手机号码长度在11到16位"
$ws.Range("F2").Value = "0.女, 1.男, 2.问号，默认为2"
$ws.Range("G2").Value = "选填，学号长度在20位之内"
$ws.Range("H2").Value = "选填，0或1，1表示禁止使用第三方登录"
$ws.Range("I2").Value = "选填，0或1，1表示禁止更改邮箱"

# --- Selection moves to K12 ---------------------------------------------
$ws.Range("K12").Select() | Out-Null
